# Auto-applies the numeric cell updates captured in the commit diff.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) receives refreshed
# market-price / profit figures (columns H-N) for specific leve rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2629.4
$ws.Range("I2").Value = 3533
$ws.Range("J2").Value = 1274
$ws.Range("K2").Value = 3533
$ws.Range("L2").Value = 1274
$ws.Range("M2").Value = -3420
$ws.Range("N2").Value = -1500

$ws.Range("H38").Value = 1217.7142
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H41").Value = 1901.4286
$ws.Range("I41").Value = 1990.625
$ws.Range("K41").Value = 1990.625
$ws.Range("M41").Value = -1550.625

$ws.Range("H58").Value = 2441
$ws.Range("J58").Value = 5500
$ws.Range("L58").Value = 16500
$ws.Range("N58").Value = -16800

$ws.Range("H61").Value = 1655
$ws.Range("I61").Value = 568.75
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 1706.25
$ws.Range("L61").Value = 18000
$ws.Range("M61").Value = -1534.25
$ws.Range("N61").Value = -18344

$ws.Range("H64").Value = 4435.8
$ws.Range("I64").Value = 3949
$ws.Range("J64").Value = 4478.1304
$ws.Range("K64").Value = 3949
$ws.Range("L64").Value = 4478.1304
$ws.Range("M64").Value = -3701
$ws.Range("N64").Value = -4974.1304

$ws.Range("H67").Value = 4435.8
$ws.Range("I67").Value = 3949
$ws.Range("J67").Value = 4478.1304
$ws.Range("K67").Value = 3949
$ws.Range("L67").Value = 4478.1304
$ws.Range("M67").Value = -3091
$ws.Range("N67").Value = -6194.1304

$ws.Range("H69").Value = 24399.2
$ws.Range("I69").Value = 15332.833
$ws.Range("K69").Value = 45998.499
$ws.Range("M69").Value = -45124.499

$ws.Range("H72").Value = 24399.2
$ws.Range("I72").Value = 15332.833
$ws.Range("K72").Value = 137995.497
$ws.Range("M72").Value = -133627.497

$ws.Range("H76").Value = 5860.1
$ws.Range("I76").Value = 3516
$ws.Range("J76").Value = 11329.667
$ws.Range("K76").Value = 3516
$ws.Range("L76").Value = 11329.667
$ws.Range("M76").Value = -3201
$ws.Range("N76").Value = -11959.667

$ws.Range("H79").Value = 5860.1
$ws.Range("I79").Value = 3516
$ws.Range("J79").Value = 11329.667
$ws.Range("K79").Value = 3516
$ws.Range("L79").Value = 11329.667
$ws.Range("M79").Value = -2424
$ws.Range("N79").Value = -13513.667

$ws.Range("H86").Value = 2520901
$ws.Range("J86").Value = 12075
$ws.Range("L86").Value = 12075
$ws.Range("N86").Value = -14321

$ws.Range("H87").Value = 154000
$ws.Range("J87").Value = 154000
$ws.Range("L87").Value = 154000
$ws.Range("N87").Value = -156496

$ws.Range("H89").Value = 2520901
$ws.Range("J89").Value = 12075
$ws.Range("L89").Value = 60375
$ws.Range("N89").Value = -71607

$ws.Range("H90").Value = 154000
$ws.Range("J90").Value = 154000
$ws.Range("L90").Value = 462000
$ws.Range("N90").Value = -474480

$ws.Range("H138").Value = 2971.7253
$ws.Range("J138").Value = 3224.9211
$ws.Range("L138").Value = 9674.763300000001
$ws.Range("N138").Value = -19954.7633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2229
$ws.Range("I2").Value = 1995.4
$ws.Range("J2").Value = 2813
$ws.Range("K2").Value = 1995.4
$ws.Range("L2").Value = 2813
$ws.Range("M2").Value = -1882.4
$ws.Range("N2").Value = -3039

$ws.Range("H24").Value = 121500
$ws.Range("J24").Value = 121500
$ws.Range("L24").Value = 121500
$ws.Range("N24").Value = -122248

$ws.Range("H31").Value = 32506.125
$ws.Range("I31").Value = 6009.8335
$ws.Range("J31").Value = 111995
$ws.Range("K31").Value = 6009.8335
$ws.Range("L31").Value = 111995
$ws.Range("M31").Value = -5715.8335
$ws.Range("N31").Value = -112583

$ws.Range("H45").Value = 2656.353
$ws.Range("I45").Value = 2410.5334
$ws.Range("J45").Value = 4500
$ws.Range("K45").Value = 2410.5334
$ws.Range("L45").Value = 4500
$ws.Range("M45").Value = -2033.5334
$ws.Range("N45").Value = -5254

$ws.Range("H61").Value = 100207600
$ws.Range("I61").Value = 125009496
$ws.Range("K61").Value = 125009496
$ws.Range("M61").Value = -125009284

$ws.Range("H74").Value = 11914856
$ws.Range("I74").Value = 16668819
$ws.Range("J74").Value = 29947.5
$ws.Range("K74").Value = 16668819
$ws.Range("L74").Value = 29947.5
$ws.Range("M74").Value = -16667945
$ws.Range("N74").Value = -31695.5

$ws.Range("H77").Value = 11914856
$ws.Range("I77").Value = 16668819
$ws.Range("J77").Value = 29947.5
$ws.Range("K77").Value = 83344095
$ws.Range("L77").Value = 149737.5
$ws.Range("M77").Value = -83339727
$ws.Range("N77").Value = -158473.5

$ws.Range("H100").Value = 121500
$ws.Range("J100").Value = 121500
$ws.Range("L100").Value = 121500
$ws.Range("N100").Value = -123664

$ws.Range("H110").Value = 1757.7273
$ws.Range("I110").Value = 1950.7778
$ws.Range("K110").Value = 1950.7778
$ws.Range("M110").Value = 94.22219999999993

$ws.Range("H116").Value = 2229
$ws.Range("I116").Value = 1995.4
$ws.Range("J116").Value = 2813
$ws.Range("K116").Value = 1995.4
$ws.Range("L116").Value = 2813
$ws.Range("M116").Value = 298.5999999999999
$ws.Range("N116").Value = -7401

$ws.Range("H136").Value = 100207600
$ws.Range("I136").Value = 125009496
$ws.Range("K136").Value = 375028488
$ws.Range("M136").Value = -375025938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2229
$ws.Range("I3").Value = 1995.4
$ws.Range("J3").Value = 2813
$ws.Range("K3").Value = 1995.4
$ws.Range("L3").Value = 2813
$ws.Range("M3").Value = -1881.4
$ws.Range("N3").Value = -3041

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H102").Value = 80355.836
$ws.Range("I102").Value = 41278
$ws.Range("K102").Value = 41278
$ws.Range("M102").Value = -38033

$ws.Range("H134").Value = 74117.64
$ws.Range("I134").Value = 2269.0833
$ws.Range("K134").Value = 6807.249899999999
$ws.Range("M134").Value = -4272.249899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 24001.8
$ws.Range("I10").Value = 5333.6665
$ws.Range("J10").Value = 52004
$ws.Range("K10").Value = 5333.6665
$ws.Range("L10").Value = 52004
$ws.Range("M10").Value = -5194.6665
$ws.Range("N10").Value = -52282

$ws.Range("H31").Value = 738440.75
$ws.Range("J31").Value = 1117061.2
$ws.Range("L31").Value = 1117061.2
$ws.Range("N31").Value = -1117651.2

$ws.Range("H34").Value = 738440.75
$ws.Range("J34").Value = 1117061.2
$ws.Range("L34").Value = 1117061.2
$ws.Range("N34").Value = -1117465.2

$ws.Range("H132").Value = 4160.4116
$ws.Range("I132").Value = 2823.2
$ws.Range("K132").Value = 8469.599999999999
$ws.Range("M132").Value = -5939.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1562.9231
$ws.Range("I113").Value = 798.5
$ws.Range("J113").Value = 1701.909
$ws.Range("K113").Value = 2395.5
$ws.Range("L113").Value = 5105.727000000001
$ws.Range("M113").Value = -225.5
$ws.Range("N113").Value = -9445.727000000001

$ws.Range("H131").Value = 3219.6667
$ws.Range("J131").Value = 3579.9048
$ws.Range("L131").Value = 10739.7144
$ws.Range("N131").Value = -20819.7144

$ws.Range("H133").Value = 5585.7144
$ws.Range("I133").Value = 4033.3333
$ws.Range("K133").Value = 12099.9999
$ws.Range("M133").Value = -7039.999899999999

$ws.Range("H140").Value = 2706.9375
$ws.Range("J140").Value = 3000
$ws.Range("L140").Value = 9000
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 18040.8
$ws.Range("J99").Value = 27801.334
$ws.Range("L99").Value = 27801.334
$ws.Range("N99").Value = -32293.334

$ws.Range("H122").Value = 1375.6428
$ws.Range("I122").Value = 1112.3
$ws.Range("K122").Value = 3336.9
$ws.Range("M122").Value = -886.8999999999996

$ws.Range("H128").Value = 62635.715
$ws.Range("J128").Value = 62635.715
$ws.Range("L128").Value = 62635.715
$ws.Range("N128").Value = -72595.715

$ws.Range("H141").Value = 24749.75
$ws.Range("J141").Value = 24749.75
$ws.Range("L141").Value = 24749.75
$ws.Range("N141").Value = -35109.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2031
$ws.Range("I16").Value = 1879.5385
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 1879.5385
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -1709.5385
$ws.Range("N16").Value = -4340

$ws.Range("H36").Value = 73808.336
$ws.Range("J36").Value = 73808.336
$ws.Range("L36").Value = 73808.336
$ws.Range("N36").Value = -74932.336

$ws.Range("H40").Value = 4237.5415
$ws.Range("I40").Value = 3778.7368
$ws.Range("K40").Value = 3778.7368
$ws.Range("M40").Value = -3642.7368

$ws.Range("H132").Value = 60929.26
$ws.Range("I132").Value = 10742.385
$ws.Range("K132").Value = 32227.155
$ws.Range("M132").Value = -29697.155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 16053.235
$ws.Range("I122").Value = 12425.137
$ws.Range("K122").Value = 37275.411
$ws.Range("M122").Value = -34825.411

$ws.Range("H140").Value = 59714.5
$ws.Range("J140").Value = 59714.5
$ws.Range("L140").Value = 59714.5
$ws.Range("N140").Value = -70074.5

$ws.Range("H141").Value = 64999.75
$ws.Range("J141").Value = 64999.75
$ws.Range("L141").Value = 64999.75
$ws.Range("N141").Value = -75359.75
